$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap F:V data between paired rows (both teams recorded same match twice; ordering corrected) ---
# Row 17 <-> Row 18
$ws.Cells.Item(17,6).Value = "FC Bhayangkara"
$ws.Cells.Item(17,7).Value = 1
$ws.Cells.Item(17,8).Value = "RANS Nusantara"
$ws.Cells.Item(17,9).Value = 2
$ws.Cells.Item(17,10).Value = 1.32
$ws.Cells.Item(17,11).Value = "08/07/2023 02:12"
$ws.Cells.Item(17,12).Value = 1.75
$ws.Cells.Item(17,13).Value = "09/07/2023 13:54"
$ws.Cells.Item(17,14).Value = 5.37
$ws.Cells.Item(17,15).Value = "08/07/2023 02:12"
$ws.Cells.Item(17,16).Value = 3.97
$ws.Cells.Item(17,17).Value = "09/07/2023 13:54"
$ws.Cells.Item(17,18).Value = 6.1
$ws.Cells.Item(17,19).Value = "08/07/2023 02:12"
$ws.Cells.Item(17,20).Value = 4.29
$ws.Cells.Item(17,21).Value = "09/07/2023 13:54"
$ws.Cells.Item(17,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/fc-bhayangkara-rans-nusantara/6Nzk5BLE/"
$ws.Cells.Item(18,6).Value = "Persikabo 1973"
$ws.Cells.Item(18,7).Value = 0
$ws.Cells.Item(18,8).Value = "Persija Jakarta"
$ws.Cells.Item(18,9).Value = 0
$ws.Cells.Item(18,10).Value = 3.34
$ws.Cells.Item(18,11).Value = "08/07/2023 02:12"
$ws.Cells.Item(18,12).Value = 4.55
$ws.Cells.Item(18,13).Value = "09/07/2023 13:59"
$ws.Cells.Item(18,14).Value = 3.3
$ws.Cells.Item(18,15).Value = "08/07/2023 02:12"
$ws.Cells.Item(18,16).Value = 3.38
$ws.Cells.Item(18,17).Value = "09/07/2023 13:59"
$ws.Cells.Item(18,18).Value = 1.99
$ws.Cells.Item(18,19).Value = "08/07/2023 02:12"
$ws.Cells.Item(18,20).Value = 1.85
$ws.Cells.Item(18,21).Value = "09/07/2023 13:59"
$ws.Cells.Item(18,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/persikabo-1973-persija-jakarta/SWyo6i68/"

# Row 19 <-> Row 20
$ws.Cells.Item(19,6).Value = "Persikabo 1973"
$ws.Cells.Item(19,7).Value = 0
$ws.Cells.Item(19,8).Value = "PSM Makassar"
$ws.Cells.Item(19,9).Value = 1
$ws.Cells.Item(19,10).Value = 3.73
$ws.Cells.Item(19,11).Value = "12/07/2023 22:12"
$ws.Cells.Item(19,12).Value = 3.53
$ws.Cells.Item(19,13).Value = "14/07/2023 09:53"
$ws.Cells.Item(19,14).Value = 3.49
$ws.Cells.Item(19,15).Value = "12/07/2023 22:12"
$ws.Cells.Item(19,16).Value = 3.33
$ws.Cells.Item(19,17).Value = "14/07/2023 09:53"
$ws.Cells.Item(19,18).Value = 1.86
$ws.Cells.Item(19,19).Value = "12/07/2023 22:12"
$ws.Cells.Item(19,20).Value = 2.11
$ws.Cells.Item(19,21).Value = "14/07/2023 09:53"
$ws.Cells.Item(19,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/persikabo-1973-psm-makassar/QPjj3kjR/"
$ws.Cells.Item(20,6).Value = "Barito Putera"
$ws.Cells.Item(20,7).Value = 3
$ws.Cells.Item(20,8).Value = "PSS Sleman"
$ws.Cells.Item(20,9).Value = 1
$ws.Cells.Item(20,10).Value = 1.56
$ws.Cells.Item(20,11).Value = "12/07/2023 22:12"
$ws.Cells.Item(20,12).Value = 1.86
$ws.Cells.Item(20,13).Value = "14/07/2023 09:58"
$ws.Cells.Item(20,14).Value = 4.07
$ws.Cells.Item(20,15).Value = "12/07/2023 22:12"
$ws.Cells.Item(20,16).Value = 3.27
$ws.Cells.Item(20,17).Value = "14/07/2023 09:57"
$ws.Cells.Item(20,18).Value = 4.9
$ws.Cells.Item(20,19).Value = "12/07/2023 22:12"
$ws.Cells.Item(20,20).Value = 3.05
$ws.Cells.Item(20,21).Value = "14/07/2023 09:58"
$ws.Cells.Item(20,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/ps-barito-putera-pss-sleman/nZin4VzL/"

# Row 22 <-> Row 23
$ws.Cells.Item(22,6).Value = "RANS Nusantara"
$ws.Cells.Item(22,7).Value = 0
$ws.Cells.Item(22,8).Value = "Persita"
$ws.Cells.Item(22,9).Value = 1
$ws.Cells.Item(22,10).Value = 2.52
$ws.Cells.Item(22,11).Value = "13/07/2023 22:12"
$ws.Cells.Item(22,12).Value = 2.85
$ws.Cells.Item(22,13).Value = "15/07/2023 09:59"
$ws.Cells.Item(22,14).Value = 3.33
$ws.Cells.Item(22,15).Value = "13/07/2023 22:12"
$ws.Cells.Item(22,16).Value = 3.45
$ws.Cells.Item(22,17).Value = "15/07/2023 09:50"
$ws.Cells.Item(22,18).Value = 2.54
$ws.Cells.Item(22,19).Value = "13/07/2023 22:12"
$ws.Cells.Item(22,20).Value = 2.4
$ws.Cells.Item(22,21).Value = "15/07/2023 09:59"
$ws.Cells.Item(22,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/rans-nusantara-persita/nczZfRc7/"
$ws.Cells.Item(23,6).Value = "Persik Kediri"
$ws.Cells.Item(23,7).Value = 5
$ws.Cells.Item(23,8).Value = "Arema FC"
$ws.Cells.Item(23,9).Value = 2
$ws.Cells.Item(23,10).Value = 1.71
$ws.Cells.Item(23,11).Value = "13/07/2023 22:12"
$ws.Cells.Item(23,12).Value = 2.08
$ws.Cells.Item(23,13).Value = "15/07/2023 09:52"
$ws.Cells.Item(23,14).Value = 3.61
$ws.Cells.Item(23,15).Value = "13/07/2023 22:12"
$ws.Cells.Item(23,16).Value = 3.29
$ws.Cells.Item(23,17).Value = "15/07/2023 09:52"
$ws.Cells.Item(23,18).Value = 4.32
$ws.Cells.Item(23,19).Value = "13/07/2023 22:12"
$ws.Cells.Item(23,20).Value = 3.66
$ws.Cells.Item(23,21).Value = "15/07/2023 09:52"
$ws.Cells.Item(23,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/persik-kediri-arema-fc/vXWQdmTf/"

# Row 42 <-> Row 43
$ws.Cells.Item(42,6).Value = "RANS Nusantara"
$ws.Cells.Item(42,7).Value = 0
$ws.Cells.Item(42,8).Value = "PSS Sleman"
$ws.Cells.Item(42,9).Value = 0
$ws.Cells.Item(42,10).Value = 2.32
$ws.Cells.Item(42,11).Value = "28/07/2023 22:12"
$ws.Cells.Item(42,12).Value = 2.79
$ws.Cells.Item(42,13).Value = "30/07/2023 09:55"
$ws.Cells.Item(42,14).Value = 3.26
$ws.Cells.Item(42,15).Value = "28/07/2023 22:12"
$ws.Cells.Item(42,16).Value = 3.5
$ws.Cells.Item(42,17).Value = "30/07/2023 09:58"
$ws.Cells.Item(42,18).Value = 2.75
$ws.Cells.Item(42,19).Value = "28/07/2023 22:12"
$ws.Cells.Item(42,20).Value = 2.42
$ws.Cells.Item(42,21).Value = "30/07/2023 09:55"
$ws.Cells.Item(42,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/rans-nusantara-pss-sleman/pCUpNqs0/"
$ws.Cells.Item(43,6).Value = "Persis Solo"
$ws.Cells.Item(43,7).Value = 1
$ws.Cells.Item(43,8).Value = "Arema FC"
$ws.Cells.Item(43,9).Value = 1
$ws.Cells.Item(43,10).Value = 1.65
$ws.Cells.Item(43,11).Value = "28/07/2023 22:12"
$ws.Cells.Item(43,12).Value = 1.61
$ws.Cells.Item(43,13).Value = "30/07/2023 09:55"
$ws.Cells.Item(43,14).Value = 3.8
$ws.Cells.Item(43,15).Value = "28/07/2023 22:12"
$ws.Cells.Item(43,16).Value = 4.15
$ws.Cells.Item(43,17).Value = "30/07/2023 09:55"
$ws.Cells.Item(43,18).Value = 4.48
$ws.Cells.Item(43,19).Value = "28/07/2023 22:12"
$ws.Cells.Item(43,20).Value = 5.06
$ws.Cells.Item(43,21).Value = "30/07/2023 09:51"
$ws.Cells.Item(43,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/persis-solo-arema-fc/tvOgLNBC/"

# Row 44 <-> Row 45
$ws.Cells.Item(44,6).Value = "Persija Jakarta"
$ws.Cells.Item(44,7).Value = 1
$ws.Cells.Item(44,8).Value = "Persebaya"
$ws.Cells.Item(44,9).Value = 0
$ws.Cells.Item(44,10).Value = 1.62
$ws.Cells.Item(44,11).Value = "29/07/2023 02:12"
$ws.Cells.Item(44,12).Value = 1.67
$ws.Cells.Item(44,13).Value = "30/07/2023 13:52"
$ws.Cells.Item(44,14).Value = 3.93
$ws.Cells.Item(44,15).Value = "29/07/2023 02:12"
$ws.Cells.Item(44,16).Value = 3.9
$ws.Cells.Item(44,17).Value = "30/07/2023 13:52"
$ws.Cells.Item(44,18).Value = 4.54
$ws.Cells.Item(44,19).Value = "29/07/2023 02:12"
$ws.Cells.Item(44,20).Value = 4.97
$ws.Cells.Item(44,21).Value = "30/07/2023 13:52"
$ws.Cells.Item(44,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/persija-jakarta-persebaya/2mPkM3d6/"
$ws.Cells.Item(45,6).Value = "Barito Putera"
$ws.Cells.Item(45,7).Value = 1
$ws.Cells.Item(45,8).Value = "Madura United"
$ws.Cells.Item(45,9).Value = 2
$ws.Cells.Item(45,10).Value = 2.12
$ws.Cells.Item(45,11).Value = "29/07/2023 02:12"
$ws.Cells.Item(45,12).Value = 2.34
$ws.Cells.Item(45,13).Value = "30/07/2023 13:51"
$ws.Cells.Item(45,14).Value = 3.31
$ws.Cells.Item(45,15).Value = "29/07/2023 02:12"
$ws.Cells.Item(45,16).Value = 3.44
$ws.Cells.Item(45,17).Value = "30/07/2023 13:51"
$ws.Cells.Item(45,18).Value = 3.13
$ws.Cells.Item(45,19).Value = "29/07/2023 02:12"
$ws.Cells.Item(45,20).Value = 2.94
$ws.Cells.Item(45,21).Value = "30/07/2023 13:51"
$ws.Cells.Item(45,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/ps-barito-putera-madura-united/z1NcKsRI/"

# Row 84 <-> Row 85
$ws.Cells.Item(84,6).Value = "Borneo"
$ws.Cells.Item(84,7).Value = 2
$ws.Cells.Item(84,8).Value = "Persita"
$ws.Cells.Item(84,9).Value = 1
$ws.Cells.Item(84,10).Value = 1.85
$ws.Cells.Item(84,11).Value = "24/08/2023 02:12"
$ws.Cells.Item(84,12).Value = 1.64
$ws.Cells.Item(84,13).Value = "25/08/2023 13:51"
$ws.Cells.Item(84,14).Value = 3.73
$ws.Cells.Item(84,15).Value = "24/08/2023 02:12"
$ws.Cells.Item(84,16).Value = 3.9
$ws.Cells.Item(84,17).Value = "25/08/2023 13:51"
$ws.Cells.Item(84,18).Value = 3.51
$ws.Cells.Item(84,19).Value = "24/08/2023 02:12"
$ws.Cells.Item(84,20).Value = 5.27
$ws.Cells.Item(84,21).Value = "25/08/2023 13:51"
$ws.Cells.Item(84,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/borneo-persita/l6QcwHle/"
$ws.Cells.Item(85,6).Value = "Dewa United"
$ws.Cells.Item(85,7).Value = 2
$ws.Cells.Item(85,8).Value = "Persija Jakarta"
$ws.Cells.Item(85,9).Value = 0
$ws.Cells.Item(85,10).Value = 3.11
$ws.Cells.Item(85,11).Value = "24/08/2023 02:12"
$ws.Cells.Item(85,12).Value = 2.6
$ws.Cells.Item(85,13).Value = "25/08/2023 13:59"
$ws.Cells.Item(85,14).Value = 3.17
$ws.Cells.Item(85,15).Value = "24/08/2023 02:12"
$ws.Cells.Item(85,16).Value = 3.11
$ws.Cells.Item(85,17).Value = "25/08/2023 13:59"
$ws.Cells.Item(85,18).Value = 2.15
$ws.Cells.Item(85,19).Value = "24/08/2023 02:12"
$ws.Cells.Item(85,20).Value = 2.84
$ws.Cells.Item(85,21).Value = "25/08/2023 13:59"
$ws.Cells.Item(85,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/dewa-united-persija-jakarta/IuV1xy41/"

# Row 101 <-> Row 102
$ws.Cells.Item(101,6).Value = "Dewa United"
$ws.Cells.Item(101,7).Value = 2
$ws.Cells.Item(101,8).Value = "FC Bhayangkara"
$ws.Cells.Item(101,9).Value = 2
$ws.Cells.Item(101,10).Value = 1.95
$ws.Cells.Item(101,11).Value = "14/09/2023 02:12"
$ws.Cells.Item(101,12).Value = 1.76
$ws.Cells.Item(101,13).Value = "15/09/2023 13:59"
$ws.Cells.Item(101,14).Value = 3.36
$ws.Cells.Item(101,15).Value = "14/09/2023 02:12"
$ws.Cells.Item(101,16).Value = 3.62
$ws.Cells.Item(101,17).Value = "15/09/2023 13:59"
$ws.Cells.Item(101,18).Value = 3.41
$ws.Cells.Item(101,19).Value = "14/09/2023 02:12"
$ws.Cells.Item(101,20).Value = 4.68
$ws.Cells.Item(101,21).Value = "15/09/2023 13:58"
$ws.Cells.Item(101,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/dewa-united-fc-bhayangkara/d4T3hmfa/"
$ws.Cells.Item(102,6).Value = "Bali United"
$ws.Cells.Item(102,7).Value = 1
$ws.Cells.Item(102,8).Value = "RANS Nusantara"
$ws.Cells.Item(102,9).Value = 2
$ws.Cells.Item(102,10).Value = 1.56
$ws.Cells.Item(102,11).Value = "14/09/2023 02:12"
$ws.Cells.Item(102,12).Value = 1.75
$ws.Cells.Item(102,13).Value = "15/09/2023 13:58"
$ws.Cells.Item(102,14).Value = 4.45
$ws.Cells.Item(102,15).Value = "14/09/2023 02:12"
$ws.Cells.Item(102,16).Value = 3.65
$ws.Cells.Item(102,17).Value = "15/09/2023 13:58"
$ws.Cells.Item(102,18).Value = 4.38
$ws.Cells.Item(102,19).Value = "14/09/2023 02:12"
$ws.Cells.Item(102,20).Value = 4.71
$ws.Cells.Item(102,21).Value = "15/09/2023 13:58"
$ws.Cells.Item(102,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/bali-united-rans-nusantara/zeS7i795/"

# Row 127 <-> Row 128
$ws.Cells.Item(127,6).Value = "Persikabo 1973"
$ws.Cells.Item(127,7).Value = 2
$ws.Cells.Item(127,8).Value = "Persis Solo"
$ws.Cells.Item(127,9).Value = 2
$ws.Cells.Item(127,10).Value = 2.54
$ws.Cells.Item(127,11).Value = "04/10/2023 21:12"
$ws.Cells.Item(127,12).Value = 3.74
$ws.Cells.Item(127,13).Value = "06/10/2023 09:56"
$ws.Cells.Item(127,14).Value = 3.24
$ws.Cells.Item(127,15).Value = "04/10/2023 21:12"
$ws.Cells.Item(127,16).Value = 3.73
$ws.Cells.Item(127,17).Value = "06/10/2023 09:58"
$ws.Cells.Item(127,18).Value = 2.52
$ws.Cells.Item(127,19).Value = "04/10/2023 21:12"
$ws.Cells.Item(127,20).Value = 1.92
$ws.Cells.Item(127,21).Value = "06/10/2023 09:56"
$ws.Cells.Item(127,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/persikabo-1973-persis-solo/OE3fW2x4/"
$ws.Cells.Item(128,6).Value = "RANS Nusantara"
$ws.Cells.Item(128,7).Value = 2
$ws.Cells.Item(128,8).Value = "PSIS Semarang"
$ws.Cells.Item(128,9).Value = 1
$ws.Cells.Item(128,10).Value = 2.54
$ws.Cells.Item(128,11).Value = "04/10/2023 21:12"
$ws.Cells.Item(128,12).Value = 3.07
$ws.Cells.Item(128,13).Value = "06/10/2023 09:52"
$ws.Cells.Item(128,14).Value = 3.19
$ws.Cells.Item(128,15).Value = "04/10/2023 21:12"
$ws.Cells.Item(128,16).Value = 3.19
$ws.Cells.Item(128,17).Value = "06/10/2023 09:52"
$ws.Cells.Item(128,18).Value = 2.54
$ws.Cells.Item(128,19).Value = "04/10/2023 21:12"
$ws.Cells.Item(128,20).Value = 2.39
$ws.Cells.Item(128,21).Value = "06/10/2023 09:52"
$ws.Cells.Item(128,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/rans-nusantara-psis-semarang/j15nYO7i/"

# Row 142 <-> Row 143
$ws.Cells.Item(142,6).Value = "Persita"
$ws.Cells.Item(142,7).Value = 2
$ws.Cells.Item(142,8).Value = "Persis Solo"
$ws.Cells.Item(142,9).Value = 1
$ws.Cells.Item(142,10).Value = 2.08
$ws.Cells.Item(142,11).Value = "20/10/2023 21:12"
$ws.Cells.Item(142,12).Value = 3.01
$ws.Cells.Item(142,13).Value = "22/10/2023 09:57"
$ws.Cells.Item(142,14).Value = 3.33
$ws.Cells.Item(142,15).Value = "20/10/2023 21:12"
$ws.Cells.Item(142,16).Value = 3.46
$ws.Cells.Item(142,17).Value = "22/10/2023 09:57"
$ws.Cells.Item(142,18).Value = 3.12
$ws.Cells.Item(142,19).Value = "20/10/2023 21:12"
$ws.Cells.Item(142,20).Value = 2.29
$ws.Cells.Item(142,21).Value = "22/10/2023 09:57"
$ws.Cells.Item(142,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/persita-persis-solo/IZmEPiLl/"
$ws.Cells.Item(143,6).Value = "Madura United"
$ws.Cells.Item(143,7).Value = 1
$ws.Cells.Item(143,8).Value = "Dewa United"
$ws.Cells.Item(143,9).Value = 4
$ws.Cells.Item(143,10).Value = 1.77
$ws.Cells.Item(143,11).Value = "20/10/2023 21:12"
$ws.Cells.Item(143,12).Value = 1.74
$ws.Cells.Item(143,13).Value = "22/10/2023 09:58"
$ws.Cells.Item(143,14).Value = 3.7
$ws.Cells.Item(143,15).Value = "20/10/2023 21:12"
$ws.Cells.Item(143,16).Value = 3.76
$ws.Cells.Item(143,17).Value = "22/10/2023 09:58"
$ws.Cells.Item(143,18).Value = 3.89
$ws.Cells.Item(143,19).Value = "20/10/2023 21:12"
$ws.Cells.Item(143,20).Value = 4.6
$ws.Cells.Item(143,21).Value = "22/10/2023 09:58"
$ws.Cells.Item(143,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/madura-united-dewa-united/0bwJOBze/"

# Row 156 <-> Row 157
$ws.Cells.Item(156,6).Value = "FC Bhayangkara"
$ws.Cells.Item(156,7).Value = 1
$ws.Cells.Item(156,8).Value = "PSIS Semarang"
$ws.Cells.Item(156,9).Value = 1
$ws.Cells.Item(156,10).Value = 3.01
$ws.Cells.Item(156,11).Value = "01/11/2023 01:12"
$ws.Cells.Item(156,12).Value = 3.74
$ws.Cells.Item(156,13).Value = "02/11/2023 12:56"
$ws.Cells.Item(156,14).Value = 3.27
$ws.Cells.Item(156,15).Value = "01/11/2023 01:12"
$ws.Cells.Item(156,16).Value = 3.53
$ws.Cells.Item(156,17).Value = "02/11/2023 12:59"
$ws.Cells.Item(156,18).Value = 2.16
$ws.Cells.Item(156,19).Value = "01/11/2023 01:12"
$ws.Cells.Item(156,20).Value = 1.97
$ws.Cells.Item(156,21).Value = "02/11/2023 12:58"
$ws.Cells.Item(156,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/fc-bhayangkara-psis-semarang/29YnRnCa/"
$ws.Cells.Item(157,6).Value = "Borneo"
$ws.Cells.Item(157,7).Value = 3
$ws.Cells.Item(157,8).Value = "Persik Kediri"
$ws.Cells.Item(157,9).Value = 0
$ws.Cells.Item(157,10).Value = 1.55
$ws.Cells.Item(157,11).Value = "01/11/2023 01:12"
$ws.Cells.Item(157,12).Value = 1.48
$ws.Cells.Item(157,13).Value = "02/11/2023 12:54"
$ws.Cells.Item(157,14).Value = 3.94
$ws.Cells.Item(157,15).Value = "01/11/2023 01:12"
$ws.Cells.Item(157,16).Value = 4.28
$ws.Cells.Item(157,17).Value = "02/11/2023 12:59"
$ws.Cells.Item(157,18).Value = 4.88
$ws.Cells.Item(157,19).Value = "01/11/2023 01:12"
$ws.Cells.Item(157,20).Value = 6.89
$ws.Cells.Item(157,21).Value = "02/11/2023 12:59"
$ws.Cells.Item(157,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/borneo-persik-kediri/ABg8YUJP/"

# --- Append new rows 169 and 170 (copy formatting from last existing row, then set values) ---
$ws.Range("A168:V168").Copy($ws.Range("A169:V169"))
$ws.Cells.Item(169,1).Value = 168
$ws.Cells.Item(169,2).Value = "indonesia"
$ws.Cells.Item(169,3).Value = "liga-1"
$ws.Cells.Item(169,4).Value = "2023-2024"
$ws.Cells.Item(169,5).Value = 45242.54166666666
$ws.Cells.Item(169,6).Value = "Bali United"
$ws.Cells.Item(169,7).Value = 1
$ws.Cells.Item(169,8).Value = "Borneo"
$ws.Cells.Item(169,9).Value = 2
$ws.Cells.Item(169,10).Value = 2.22
$ws.Cells.Item(169,11).Value = "11/11/2023 01:13"
$ws.Cells.Item(169,12).Value = 2.49
$ws.Cells.Item(169,13).Value = "12/11/2023 12:54"
$ws.Cells.Item(169,14).Value = 3.32
$ws.Cells.Item(169,15).Value = "11/11/2023 01:13"
$ws.Cells.Item(169,16).Value = 3.43
$ws.Cells.Item(169,17).Value = "12/11/2023 12:32"
$ws.Cells.Item(169,18).Value = 2.85
$ws.Cells.Item(169,19).Value = "11/11/2023 01:13"
$ws.Cells.Item(169,20).Value = 2.74
$ws.Cells.Item(169,21).Value = "12/11/2023 12:54"
$ws.Cells.Item(169,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/bali-united-borneo/OvBJv6ta/"

$ws.Range("A169:V169").Copy($ws.Range("A170:V170"))
$ws.Cells.Item(170,1).Value = 169
$ws.Cells.Item(170,2).Value = "indonesia"
$ws.Cells.Item(170,3).Value = "liga-1"
$ws.Cells.Item(170,4).Value = "2023-2024"
$ws.Cells.Item(170,5).Value = 45242.54166666666
$ws.Cells.Item(170,6).Value = "Dewa United"
$ws.Cells.Item(170,7).Value = 1
$ws.Cells.Item(170,8).Value = "PSM Makassar"
$ws.Cells.Item(170,9).Value = 1
$ws.Cells.Item(170,10).Value = 2.02
$ws.Cells.Item(170,11).Value = "11/11/2023 01:13"
$ws.Cells.Item(170,12).Value = 2.07
$ws.Cells.Item(170,13).Value = "12/11/2023 12:55"
$ws.Cells.Item(170,14).Value = 3.26
$ws.Cells.Item(170,15).Value = "11/11/2023 01:13"
$ws.Cells.Item(170,16).Value = 3.38
$ws.Cells.Item(170,17).Value = "12/11/2023 12:55"
$ws.Cells.Item(170,18).Value = 3.33
$ws.Cells.Item(170,19).Value = "11/11/2023 01:13"
$ws.Cells.Item(170,20).Value = 3.6
$ws.Cells.Item(170,21).Value = "12/11/2023 12:55"
$ws.Cells.Item(170,22).Value = "https://www.betexplorer.com/football/indonesia/liga-1/dewa-united-psm-makassar/2PMAtSBn/"

